$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC row 38
$ws_ALC.Range("H38").Value = 771.8182
$ws_ALC.Range("I38").Value = 538
$ws_ALC.Range("K38").Value = 1614
$ws_ALC.Range("M38").Value = -1242

# ALC row 55
$ws_ALC.Range("H55").Value = 606.8
$ws_ALC.Range("I55").Value = 1044
$ws_ALC.Range("J55").Value = 169.6
$ws_ALC.Range("K55").Value = 1044
$ws_ALC.Range("L55").Value = 169.6
$ws_ALC.Range("M55").Value = -830
$ws_ALC.Range("N55").Value = -597.6

# ALC row 86
$ws_ALC.Range("H86").Value = 8144.9414
$ws_ALC.Range("I86").Value = 2024.75
$ws_ALC.Range("J86").Value = 13585.111
$ws_ALC.Range("K86").Value = 2024.75
$ws_ALC.Range("L86").Value = 13585.111
$ws_ALC.Range("M86").Value = -901.75
$ws_ALC.Range("N86").Value = -15831.111

# ALC row 89
$ws_ALC.Range("H89").Value = 8144.9414
$ws_ALC.Range("I89").Value = 2024.75
$ws_ALC.Range("J89").Value = 13585.111
$ws_ALC.Range("K89").Value = 10123.75
$ws_ALC.Range("L89").Value = 67925.55500000001
$ws_ALC.Range("M89").Value = -4507.75
$ws_ALC.Range("N89").Value = -79157.55500000001

# ALC row 98
$ws_ALC.Range("H98").Value = 280.3889
$ws_ALC.Range("I98").Value = 280.3889
$ws_ALC.Range("K98").Value = 280.3889
$ws_ALC.Range("M98").Value = 1217.6111

# ALC row 122
$ws_ALC.Range("H122").Value = 280.3889
$ws_ALC.Range("I122").Value = 280.3889
$ws_ALC.Range("K122").Value = 841.1667
$ws_ALC.Range("M122").Value = 1608.8333

# ALC row 137
$ws_ALC.Range("H137").Value = 38331.816
$ws_ALC.Range("I137").Value = 992.55554
$ws_ALC.Range("J137").Value = 113010.336
$ws_ALC.Range("K137").Value = 2977.66662
$ws_ALC.Range("L137").Value = 339031.008
$ws_ALC.Range("M137").Value = -427.66662
$ws_ALC.Range("N137").Value = -344131.008

# ALC row 138
$ws_ALC.Range("H138").Value = 1678.68
$ws_ALC.Range("J138").Value = 2567.5557
$ws_ALC.Range("L138").Value = 7702.6671
$ws_ALC.Range("N138").Value = -17982.6671

# ARM row 32
$ws_ARM.Range("H32").Value = 17051.508
$ws_ARM.Range("I32").Value = 17274.016
$ws_ARM.Range("J32").Value = 12304.667
$ws_ARM.Range("K32").Value = 17274.016
$ws_ARM.Range("L32").Value = 12304.667
$ws_ARM.Range("M32").Value = -16987.016
$ws_ARM.Range("N32").Value = -12878.667

# ARM row 45
$ws_ARM.Range("H45").Value = 3985.0625
$ws_ARM.Range("I45").Value = 3946.923
$ws_ARM.Range("J45").Value = 4011.158
$ws_ARM.Range("K45").Value = 3946.923
$ws_ARM.Range("L45").Value = 4011.158
$ws_ARM.Range("M45").Value = -3569.923
$ws_ARM.Range("N45").Value = -4765.157999999999

# ARM row 74
$ws_ARM.Range("H74").Value = 41669440
$ws_ARM.Range("I74").Value = 52634324
$ws_ARM.Range("K74").Value = 52634324
$ws_ARM.Range("M74").Value = -52633450

# ARM row 77
$ws_ARM.Range("H77").Value = 41669440
$ws_ARM.Range("I77").Value = 52634324
$ws_ARM.Range("K77").Value = 263171620
$ws_ARM.Range("M77").Value = -263167252

# BSM row 86
$ws_BSM.Range("H86").Value = 1474.3959
$ws_BSM.Range("I86").Value = 1386.079
$ws_BSM.Range("K86").Value = 1386.079
$ws_BSM.Range("M86").Value = -263.079

# BSM row 89
$ws_BSM.Range("H89").Value = 1474.3959
$ws_BSM.Range("I89").Value = 1386.079
$ws_BSM.Range("K89").Value = 6930.395
$ws_BSM.Range("M89").Value = -1314.395

# CRP row 58
$ws_CRP.Range("H58").Value = 27751.422
$ws_CRP.Range("I58").Value = 1666.6154
$ws_CRP.Range("J58").Value = 84268.5
$ws_CRP.Range("K58").Value = 1666.6154
$ws_CRP.Range("L58").Value = 84268.5
$ws_CRP.Range("M58").Value = -1463.6154
$ws_CRP.Range("N58").Value = -84674.5

# CRP row 62
$ws_CRP.Range("H62").Value = 5746.5
$ws_CRP.Range("I62").Value = 3999.5
$ws_CRP.Range("J62").Value = 6620
$ws_CRP.Range("K62").Value = 3999.5
$ws_CRP.Range("L62").Value = 6620
$ws_CRP.Range("M62").Value = -3375.5
$ws_CRP.Range("N62").Value = -7868

# CRP row 65
$ws_CRP.Range("H65").Value = 5746.5
$ws_CRP.Range("I65").Value = 3999.5
$ws_CRP.Range("J65").Value = 6620
$ws_CRP.Range("K65").Value = 19997.5
$ws_CRP.Range("L65").Value = 33100
$ws_CRP.Range("M65").Value = -16877.5
$ws_CRP.Range("N65").Value = -39340

# CRP row 122
$ws_CRP.Range("H122").Value = 1792
$ws_CRP.Range("I122").Value = 1746.5
$ws_CRP.Range("K122").Value = 5239.5
$ws_CRP.Range("M122").Value = -2789.5

# CRP row 134
$ws_CRP.Range("H134").Value = 1033.862
$ws_CRP.Range("I134").Value = 879.75
$ws_CRP.Range("J134").Value = 1223.5385
$ws_CRP.Range("K134").Value = 2639.25
$ws_CRP.Range("L134").Value = 3670.6155
$ws_CRP.Range("M134").Value = -104.25
$ws_CRP.Range("N134").Value = -8740.6155

# CRP row 136
$ws_CRP.Range("H136").Value = 27751.422
$ws_CRP.Range("I136").Value = 1666.6154
$ws_CRP.Range("J136").Value = 84268.5
$ws_CRP.Range("K136").Value = 4999.8462
$ws_CRP.Range("L136").Value = 252805.5
$ws_CRP.Range("M136").Value = -2449.8462
$ws_CRP.Range("N136").Value = -257905.5

# CUL row 3
$ws_CUL.Range("H3").Value = 3821.0833
$ws_CUL.Range("J3").Value = 4713.2856
$ws_CUL.Range("L3").Value = 14139.8568
$ws_CUL.Range("N3").Value = -14363.8568

# CUL row 4
$ws_CUL.Range("H4").Value = 2500084
$ws_CUL.Range("I4").Value = 78.666664
$ws_CUL.Range("K4").Value = 235.999992
$ws_CUL.Range("M4").Value = -123.999992

# CUL row 23
$ws_CUL.Range("H23").Value = 571.8182
$ws_CUL.Range("J23").Value = 822.2308
$ws_CUL.Range("L23").Value = 2466.6924
$ws_CUL.Range("N23").Value = -2936.6924

# CUL row 116
$ws_CUL.Range("H116").Value = 1019.8
$ws_CUL.Range("I116").Value = 774.75
$ws_CUL.Range("K116").Value = 2324.25
$ws_CUL.Range("M116").Value = 1117.75

# CUL row 122
$ws_CUL.Range("H122").Value = 574
$ws_CUL.Range("I122").Value = 227.14285
$ws_CUL.Range("J122").Value = 1059.6
$ws_CUL.Range("K122").Value = 2044.28565
$ws_CUL.Range("L122").Value = 9536.4
$ws_CUL.Range("M122").Value = 405.71435
$ws_CUL.Range("N122").Value = -14436.4

# CUL row 131
$ws_CUL.Range("H131").Value = 672.1
$ws_CUL.Range("J131").Value = 675.65625
$ws_CUL.Range("L131").Value = 2026.96875
$ws_CUL.Range("N131").Value = -12106.96875

# CUL row 132
$ws_CUL.Range("H132").Value = 1078.9375
$ws_CUL.Range("J132").Value = 1308
$ws_CUL.Range("L132").Value = 11772
$ws_CUL.Range("N132").Value = -16832

# CUL row 140
$ws_CUL.Range("H140").Value = 1665.3572
$ws_CUL.Range("I140").Value = 1398.3334
$ws_CUL.Range("K140").Value = 4195.0002
$ws_CUL.Range("M140").Value = 984.9997999999996

# CUL row 141
$ws_CUL.Range("H141").Value = 3465.0833
$ws_CUL.Range("I141").Value = 5517.6
$ws_CUL.Range("K141").Value = 16552.8
$ws_CUL.Range("M141").Value = -11372.8

# GSM row 70
$ws_GSM.Range("H70").Value = 2985438.5
$ws_GSM.Range("I70").Value = 13636.818
$ws_GSM.Range("J70").Value = 6254420.5
$ws_GSM.Range("K70").Value = 13636.818
$ws_GSM.Range("L70").Value = 6254420.5
$ws_GSM.Range("M70").Value = -13366.818
$ws_GSM.Range("N70").Value = -6254960.5

# GSM row 73
$ws_GSM.Range("H73").Value = 2985438.5
$ws_GSM.Range("I73").Value = 13636.818
$ws_GSM.Range("J73").Value = 6254420.5
$ws_GSM.Range("K73").Value = 13636.818
$ws_GSM.Range("L73").Value = 6254420.5
$ws_GSM.Range("M73").Value = -12700.818
$ws_GSM.Range("N73").Value = -6256292.5

# GSM row 101
$ws_GSM.Range("H101").Value = 0
$ws_GSM.Range("J101").Value = 0
$ws_GSM.Range("N101").ClearContents()

# LTW row 7
$ws_LTW.Range("H7").Value = 4479.0713
$ws_LTW.Range("I7").Value = 4333.9165
$ws_LTW.Range("J7").Value = 5350
$ws_LTW.Range("K7").Value = 4333.9165
$ws_LTW.Range("L7").Value = 5350
$ws_LTW.Range("M7").Value = -4221.9165
$ws_LTW.Range("N7").Value = -5574

# LTW row 40
$ws_LTW.Range("H40").Value = 3572
$ws_LTW.Range("I40").Value = 1959.8
$ws_LTW.Range("J40").Value = 4646.8
$ws_LTW.Range("K40").Value = 1959.8
$ws_LTW.Range("L40").Value = 4646.8
$ws_LTW.Range("M40").Value = -1823.8
$ws_LTW.Range("N40").Value = -4918.8

# LTW row 126
$ws_LTW.Range("H126").Value = 4479.0713
$ws_LTW.Range("I126").Value = 4333.9165
$ws_LTW.Range("J126").Value = 5350
$ws_LTW.Range("K126").Value = 13001.7495
$ws_LTW.Range("L126").Value = 16050
$ws_LTW.Range("M126").Value = -10531.7495
$ws_LTW.Range("N126").Value = -20990

# LTW row 136
$ws_LTW.Range("H136").Value = 56796.777
$ws_LTW.Range("I136").Value = 72713
$ws_LTW.Range("J136").Value = 1090
$ws_LTW.Range("K136").Value = 218139
$ws_LTW.Range("L136").Value = 3270
$ws_LTW.Range("M136").Value = -215589
$ws_LTW.Range("N136").Value = -8370

# WVR row 107
$ws_WVR.Range("H107").Value = 2274030
$ws_WVR.Range("I107").Value = 985.2
$ws_WVR.Range("J107").Value = 4547074.5
$ws_WVR.Range("K107").Value = 2955.6
$ws_WVR.Range("L107").Value = 13641223.5
$ws_WVR.Range("M107").Value = -1035.6
$ws_WVR.Range("N107").Value = -13645063.5

# WVR row 136
$ws_WVR.Range("H136").Value = 41668940
$ws_WVR.Range("I136").Value = 62502124
$ws_WVR.Range("J136").Value = 2575.5
$ws_WVR.Range("K136").Value = 187506372
$ws_WVR.Range("L136").Value = 7726.5
$ws_WVR.Range("M136").Value = -187503822
$ws_WVR.Range("N136").Value = -12826.5
